$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$ws.Range("N4").Value = "nan"
$ws.Range("O4").Value = "nan"
